$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.60633566666667
$ws.Range("H2").Value = 52.819007
$ws.Range("I2").Value = 0.01967183396478544
$ws.Range("J2").Value = 0.01967183396478544
$ws.Range("M2").Value = 12.23245266666667
$ws.Range("N2").Value = 36.697358
$ws.Range("O2").Value = 0.2605176191876535
$ws.Range("P2").Value = 0.2605176191876535
$ws.Range("Q2").Value = 215.3686676759451
$ws.Range("R2").Value = 1938.318009083506
$ws.Range("S2").Value = 0.005124859349560721
$ws.Range("T2").Value = 0.00512485934956072
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.60633566666667
$ws.Range("H3").Value = 52.819007
$ws.Range("I3").Value = 0.01967183396478544
$ws.Range("J3").Value = 0.01967183396478544
$ws.Range("O3").Value = 0.2463118768921188
$ws.Range("P3").Value = 0.2463118768921188
$ws.Range("Q3").Value = 203.6248485781152
$ws.Range("R3").Value = 1832.623637203037
$ws.Range("S3").Value = 0.004845406345776432
$ws.Range("T3").Value = 0.004845406345776431
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.60633566666667
$ws.Range("H4").Value = 52.819007
$ws.Range("I4").Value = 0.01967183396478544
$ws.Range("J4").Value = 0.01967183396478544
$ws.Range("M4").Value = 20.18080466666667
$ws.Range("N4").Value = 60.54241400000001
$ws.Range("O4").Value = 0.4297956696270414
$ws.Range("P4").Value = 0.4297956696270414
$ws.Range("Q4").Value = 355.3100209847665
$ws.Range("R4").Value = 3197.790188862899
$ws.Range("S4").Value = 0.008454869051686936
$ws.Range("T4").Value = 0.008454869051686933
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.60633566666667
$ws.Range("H5").Value = 52.819007
$ws.Range("I5").Value = 0.01967183396478544
$ws.Range("J5").Value = 0.01967183396478544
$ws.Range("M5").Value = 2.975728333333333
$ws.Range("N5").Value = 8.927185
$ws.Range("O5").Value = 0.06337483429318623
$ws.Range("P5").Value = 0.06337483429318624
$ws.Range("Q5").Value = 52.39167188947722
$ws.Range("R5").Value = 471.525047005295
$ws.Range("S5").Value = 0.00124669921776135
$ws.Range("T5").Value = 0.00124669921776135
$ws.Range("I6").Value = 0.6688940299055509
$ws.Range("J6").Value = 0.6688940299055508
$ws.Range("M6").Value = 12.23245266666667
$ws.Range("N6").Value = 36.697358
$ws.Range("O6").Value = 0.2605176191876535
$ws.Range("P6").Value = 0.2605176191876535
$ws.Range("Q6").Value = 7323.100443762999
$ws.Range("R6").Value = 65907.90399386699
$ws.Range("S6").Value = 0.1742586801598292
$ws.Range("T6").Value = 0.1742586801598292
$ws.Range("I7").Value = 0.6688940299055509
$ws.Range("J7").Value = 0.6688940299055508
$ws.Range("O7").Value = 0.2463118768921188
$ws.Range("P7").Value = 0.2463118768921188
$ws.Range("S7").Value = 0.1647565439479693
$ws.Range("T7").Value = 0.1647565439479692
$ws.Range("I8").Value = 0.6688940299055509
$ws.Range("J8").Value = 0.6688940299055508
$ws.Range("M8").Value = 20.18080466666667
$ws.Range("N8").Value = 60.54241400000001
$ws.Range("O8").Value = 0.4297956696270414
$ws.Range("P8").Value = 0.4297956696270414
$ws.Range("Q8").Value = 12081.47406224402
$ws.Range("R8").Value = 108733.2665601962
$ws.Range("S8").Value = 0.2874877574927865
$ws.Range("T8").Value = 0.2874877574927864
$ws.Range("I9").Value = 0.6688940299055509
$ws.Range("J9").Value = 0.6688940299055508
$ws.Range("M9").Value = 2.975728333333333
$ws.Range("N9").Value = 8.927185
$ws.Range("O9").Value = 0.06337483429318623
$ws.Range("P9").Value = 0.06337483429318624
$ws.Range("Q9").Value = 1781.454469693823
$ws.Range("R9").Value = 16033.09022724441
$ws.Range("S9").Value = 0.04239104830496584
$ws.Range("T9").Value = 0.04239104830496584
$ws.Range("G10").Value = 274.6625416666666
$ws.Range("H10").Value = 823.987625
$ws.Range("I10").Value = 0.3068847498029997
$ws.Range("J10").Value = 0.3068847498029996
$ws.Range("M10").Value = 12.23245266666667
$ws.Range("N10").Value = 36.697358
$ws.Range("O10").Value = 0.2605176191876535
$ws.Range("P10").Value = 0.2605176191876535
$ws.Range("Q10").Value = 3359.796540243861
$ws.Range("R10").Value = 30238.16886219475
$ws.Range("S10").Value = 0.07994888438367619
$ws.Range("T10").Value = 0.07994888438367617
$ws.Range("G11").Value = 274.6625416666666
$ws.Range("H11").Value = 823.987625
$ws.Range("I11").Value = 0.3068847498029997
$ws.Range("J11").Value = 0.3068847498029996
$ws.Range("O11").Value = 0.2463118768921188
$ws.Range("P11").Value = 0.2463118768921188
$ws.Range("Q11").Value = 3176.590490822097
$ws.Range("R11").Value = 28589.31441739888
$ws.Range("S11").Value = 0.07558935871354512
$ws.Range("T11").Value = 0.07558935871354512
$ws.Range("G12").Value = 274.6625416666666
$ws.Range("H12").Value = 823.987625
$ws.Range("I12").Value = 0.3068847498029997
$ws.Range("J12").Value = 0.3068847498029996
$ws.Range("M12").Value = 20.18080466666667
$ws.Range("N12").Value = 60.54241400000001
$ws.Range("O12").Value = 0.4297956696270414
$ws.Range("P12").Value = 0.4297956696270414
$ws.Range("Q12").Value = 5542.911102625195
$ws.Range("R12").Value = 49886.19992362676
$ws.Range("S12").Value = 0.1318977365399073
$ws.Range("T12").Value = 0.1318977365399073
$ws.Range("G13").Value = 274.6625416666666
$ws.Range("H13").Value = 823.987625
$ws.Range("I13").Value = 0.3068847498029997
$ws.Range("J13").Value = 0.3068847498029996
$ws.Range("M13").Value = 2.975728333333333
$ws.Range("N13").Value = 8.927185
$ws.Range("O13").Value = 0.06337483429318623
$ws.Range("P13").Value = 0.06337483429318624
$ws.Range("Q13").Value = 817.3211073428471
$ws.Range("R13").Value = 7355.889966085625
$ws.Range("S13").Value = 0.01944877016587102
$ws.Range("T13").Value = 0.01944877016587102
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 4.071711
$ws.Range("H14").Value = 12.215133
$ws.Range("I14").Value = 0.004549386326664026
$ws.Range("J14").Value = 0.004549386326664025
$ws.Range("M14").Value = 12.23245266666667
$ws.Range("N14").Value = 36.697358
$ws.Range("O14").Value = 0.2605176191876535
$ws.Range("P14").Value = 0.2605176191876535
$ws.Range("Q14").Value = 49.807012079846
$ws.Range("R14").Value = 448.263108718614
$ws.Range("S14").Value = 0.001185195294587376
$ws.Range("T14").Value = 0.001185195294587376
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4.071711
$ws.Range("H15").Value = 12.215133
$ws.Range("I15").Value = 0.004549386326664026
$ws.Range("J15").Value = 0.004549386326664025
$ws.Range("O15").Value = 0.2463118768921188
$ws.Range("P15").Value = 0.2463118768921188
$ws.Range("Q15").Value = 47.091089907967
$ws.Range("R15").Value = 423.819809171703
$ws.Range("S15").Value = 0.001120567884827958
$ws.Range("T15").Value = 0.001120567884827958
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 4.071711
$ws.Range("H16").Value = 12.215133
$ws.Range("I16").Value = 0.004549386326664026
$ws.Range("J16").Value = 0.004549386326664025
$ws.Range("M16").Value = 20.18080466666667
$ws.Range("N16").Value = 60.54241400000001
$ws.Range("O16").Value = 0.4297956696270414
$ws.Range("P16").Value = 0.4297956696270414
$ws.Range("Q16").Value = 82.170404350118
$ws.Range("R16").Value = 739.5336391510621
$ws.Range("S16").Value = 0.001955306542660671
$ws.Range("T16").Value = 0.00195530654266067
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 4.071711
$ws.Range("H17").Value = 12.215133
$ws.Range("I17").Value = 0.004549386326664026
$ws.Range("J17").Value = 0.004549386326664025
$ws.Range("M17").Value = 2.975728333333333
$ws.Range("N17").Value = 8.927185
$ws.Range("O17").Value = 0.06337483429318623
$ws.Range("P17").Value = 0.06337483429318624
$ws.Range("Q17").Value = 12.116305787845
$ws.Range("R17").Value = 109.046752090605
$ws.Range("S17").Value = 0.0002883166045880198
$ws.Range("T17").Value = 0.0002883166045880198
